$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(4, 7).Value = "Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Hend Farid, Dr. Amal Awwad, Dr. Shimaa Ashraf"
$ws.Cells.Item(5, 7).Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad, D Wessam Atef"
$ws.Cells.Item(6, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Amany Raafat, Dr. Marina Youhanna"
$ws.Cells.Item(8, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled"
$ws.Cells.Item(10, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(15, 7).Value = "Dr. Marian Samir, Dr. Nourhan Mohammad, Dr. Nourham Mostafa, Dr. Ahmad Mostafa, Dr. Afaf Abdallah"
$ws.Cells.Item(16, 7).Value = "Dr. Rada Rabea, Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Walaa Ghanima"
$ws.Cells.Item(17, 7).Value = "Dr. Nardine, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Monica"
$ws.Cells.Item(18, 7).Value = "Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Yasmin, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Remon, Dr. Monica"
$ws.Cells.Item(19, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(21, 7).Value = "Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Hend Farid, Dr. Amal Awwad, Dr. Shimaa Ashraf"
$ws.Cells.Item(22, 7).Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad, D Wessam Atef"
$ws.Cells.Item(23, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Amany Raafat, Dr. Marina Youhanna"
$ws.Cells.Item(25, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled"
$ws.Cells.Item(27, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(32, 7).Value = "Dr. Marian Samir, Dr. Nourhan Mohammad, Dr. Nourham Mostafa, Dr. Ahmad Mostafa, Dr. Afaf Abdallah"
$ws.Cells.Item(33, 7).Value = "Dr. Rada Rabea, Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Walaa Ghanima"
$ws.Cells.Item(34, 7).Value = "Dr. Nardine, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Monica"
$ws.Cells.Item(35, 7).Value = "Dr. Aya Emad, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Cells.Item(36, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Administrator, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(37, 7).Value = "Dr. Nada Mohammad, Dr. Kerelos Zareef, Administrator"
$ws.Cells.Item(39, 7).Value = "Dr. Shimaa Ashraf, Dr. Omnia Mohammad"
$ws.Cells.Item(40, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Amany Raafat, Dr. Marina Youhanna"
$ws.Cells.Item(44, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(45, 7).Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(48, 7).Value = "Dr. Marian Samir, Dr. Aya Alaa-Eldein, Dr. Afaf Abdallah"
$ws.Cells.Item(50, 7).Value = "Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(51, 7).Value = "Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Monica"
$ws.Cells.Item(52, 7).Value = "Dr. Naema Gomaa, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Cells.Item(53, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Administrator, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(54, 7).Value = "Dr. Nada Mohammad, Dr. Kerelos Zareef, Administrator"
$ws.Cells.Item(56, 7).Value = "Dr. Shimaa Ashraf, Dr. Omnia Mohammad"
$ws.Cells.Item(57, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Amany Raafat, Dr. Marina Youhanna"
$ws.Cells.Item(61, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(62, 7).Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(65, 7).Value = "Dr. Marian Samir, Dr. Aya Alaa-Eldein, Dr. Afaf Abdallah"
$ws.Cells.Item(67, 7).Value = "Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(68, 7).Value = "Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Monica"
$ws.Cells.Item(69, 7).Value = "Dr. Naema Gomaa, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Cells.Item(70, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(72, 7).Value = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Safa Hany"
$ws.Cells.Item(73, 7).Value = "Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Hend Farid, Dr. Amal Awwad, Dr. Shimaa Ashraf"
$ws.Cells.Item(74, 7).Value = "Dr. Aya Saeed, Dr. Amal Awwad, D Wessam Atef, Dr. Sara Nabil, Dr. Omnia Mohammad"
$ws.Cells.Item(75, 7).Value = "Dr. Eman M. Elsaid, Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Maryam Ahmad"
$ws.Cells.Item(76, 7).Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef, Dr. Mariam Toma Gerges"
$ws.Cells.Item(77, 7).Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mariam Toma Gerges"
$ws.Cells.Item(78, 7).Value = "Dr. Mohammad Safwat, Dr. Mayar Ahmad Embaby, Dr. Al-Shimaa Khaled"
$ws.Cells.Item(79, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(83, 7).Value = "Dr. Marian Samir, Dr. Youstina Ibrahim, Dr. Afaf Abdallah"
$ws.Cells.Item(84, 7).Value = "Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(85, 7).Value = "Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Nahla, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Monica"
$ws.Cells.Item(86, 7).Value = "Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Yasmin, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Remon, Dr. Monica"
$ws.Cells.Item(87, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(88, 7).Value = "Dr. Nada Mohammad, Dr. Fatma Elhady"
$ws.Cells.Item(89, 7).Value = "Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Hend Farid, Dr. Amal Awwad, Dr. Shimaa Ashraf"
$ws.Cells.Item(90, 7).Value = "Dr. Aya Saeed, Dr. Amal Awwad, D Wessam Atef, Dr. Sara Nabil, Dr. Omnia Mohammad"
$ws.Cells.Item(91, 7).Value = "Dr. Eman M. Elsaid, Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Maryam Ahmad"
$ws.Cells.Item(92, 7).Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mariam Toma Gerges"
$ws.Cells.Item(94, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(98, 7).Value = "Dr. Walaa Ghanima, Dr. Nourhan Mohammad, Dr. Afaf Abdallah, Dr. Nourhan Hosni"
$ws.Cells.Item(101, 7).Value = "Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(102, 7).Value = "Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Nahla, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Monica"
$ws.Cells.Item(103, 7).Value = "Dr. Aya Hanafy, Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Gehad Salah, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Abdullah El-Agrody"
$ws.Cells.Item(104, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Administrator, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(106, 7).Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad, D Wessam Atef"
$ws.Cells.Item(107, 7).Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad"
$ws.Cells.Item(108, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Amany Raafat"
$ws.Cells.Item(110, 7).Value = "Dr. Mohammad Safwat, Dr. Afnan Fares, Dr. Al-Shimaa Khaled"
$ws.Cells.Item(112, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(116, 7).Value = "Dr. Enas Omran, Dr. Afaf Abdallah, Dr. Nourham Mostafa"
$ws.Cells.Item(117, 7).Value = "Dr. Amr Saeed, Dr. Enas Omran, Dr. Taqwa Mohammad"
$ws.Cells.Item(119, 7).Value = "Dr. Aya Hanafy, Dr. Nardine, Dr. Shorok Mohammad, Dr. Nahla, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon"
$ws.Cells.Item(120, 7).Value = "Dr. Aya Hanafy, Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Gehad Salah, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Abdullah El-Agrody"
$ws.Cells.Item(121, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Administrator, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(123, 7).Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad, D Wessam Atef"
$ws.Cells.Item(124, 7).Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad"
$ws.Cells.Item(127, 7).Value = "Dr. Mohammad Safwat, Dr. Afnan Fares, Dr. Al-Shimaa Khaled"
$ws.Cells.Item(128, 7).Value = "Dr. Mohammad Safwat, Dr. Mayar Ahmad Embaby, Dr. Al-Shimaa Khaled"
$ws.Cells.Item(129, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Al-Shimaa Khaled, Administrator"
$ws.Cells.Item(133, 7).Value = "Dr. Enas Omran, Dr. Afaf Abdallah, Dr. Nourham Mostafa"
$ws.Cells.Item(134, 7).Value = "Dr. Amr Saeed, Dr. Enas Omran, Dr. Taqwa Mohammad"
$ws.Cells.Item(136, 7).Value = "Dr. Aya Hanafy, Dr. Nardine, Dr. Shorok Mohammad, Dr. Nahla, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon"
$ws.Cells.Item(137, 7).Value = "Dr. Aya Hanafy, Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Gehad Salah, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Abdullah El-Agrody"
$ws.Cells.Item(138, 7).Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Administrator, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Cells.Item(142, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Amany Raafat"
$ws.Cells.Item(144, 7).Value = "Dr. Mohammad Safwat, Nourhan Mamdouh Hassan, Dr. Mayar Ahmad Embaby, Dr. Mariam Toma Gerges"
$ws.Cells.Item(145, 7).Value = "Dr. Mohammad Safwat, Dr. Mayar Ahmad Embaby, Dr. Al-Shimaa Khaled"
$ws.Cells.Item(150, 7).Value = "Dr. Marian Samir, Dr. Youstina Ibrahim, Dr. Afaf Abdallah"
$ws.Cells.Item(151, 7).Value = "Dr. Hana Amr, Dr. Marian Samir, Administrator, Dr. Rada Rabea, Dr. Nourhan Mohammad"
